# Initial support of 2 HX711 working together... not tested yet.
# Duplicate the "1 Sensor" calculation block (rows 1-18) into a second
# "Sensor 2" block at rows 21-38, relabel the first block's header as
# "Sensor 1", and move the active-sheet/selection state onto "1 Sensor".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1 Sensor")
$ws2 = $wb.Worksheets.Item("2 Sensors")

# ---------------------------------------------------------------------
# 1) Re-label the existing block's header: "1 Sensor" -> "Sensor 1", and
#    extend the orange "Input" header band across B1:J1 (was only A1).
# ---------------------------------------------------------------------
$ws1.Range("A1").Value = "Sensor 1"
$ws1.Range("A1:J1").Style = "Input"

# ---------------------------------------------------------------------
# 2) Drop the old stray leftover cell at A22 (not part of the new layout).
# ---------------------------------------------------------------------
$ws1.Range("A22").ClearContents()

# ---------------------------------------------------------------------
# 3) New "Sensor 2" header band at row 21 (mirrors row 1).
# ---------------------------------------------------------------------
$ws1.Range("A21").Value = "Sensor 2"
$ws1.Range("A21:J21").Style = "Input"

# ---------------------------------------------------------------------
# 4) Duplicate the "Ref weights / Offsets" table (rows 3-4) to rows 23-24.
# ---------------------------------------------------------------------
$ws1.Range("A23").Value = "Ref weights:"
$ws1.Range("B23").Value = 25
$ws1.Range("C23").Value = 75
$ws1.Range("D23").Value = 150
$ws1.Range("E23").Value = 225
$ws1.Range("F23").Value = 300
$ws1.Range("G23").Value = 375
$ws1.Range("H23").Value = 475
$ws1.Range("J23").Value = "0ffset (Sensor):"

$ws1.Range("A24").Value = "Offsets:"
$ws1.Range("B24").Value = 8670386.5332999993
$ws1.Range("C24").Value = 8752399.9066700004
$ws1.Range("D24").Value = 8876659.7200000007
$ws1.Range("E24").Value = 9001015.1866699997
$ws1.Range("F24").Value = 9125726.5800000001
$ws1.Range("G24").Value = 9249752.4866700005
$ws1.Range("H24").Value = 9415005.6466700006
$ws1.Range("J24").Value = 8628270.3399999999
$ws1.Range("B24:H24").NumberFormat = "0.000"
$ws1.Range("J24").NumberFormat = "0.000"

# ---------------------------------------------------------------------
# 5) Duplicate the "Average ratios / Ratios" table (rows 6-7) to rows 26-27.
# ---------------------------------------------------------------------
$ws1.Range("J26").Value = "Average ratios:"

$ws1.Range("A27").Value = "Ratios:"
$ws1.Range("B27").Formula = "=(B24-J24)/B23"
$ws1.Range("C27").Formula = "=(C24-J24)/C23"
$ws1.Range("D27").Formula = "=(D24-J24)/D23"
$ws1.Range("E27").Formula = "=(E24-J24)/E23"
$ws1.Range("F27").Formula = "=(F24-J24)/F23"
$ws1.Range("G27").Formula = "=(G24-J24)/G23"
$ws1.Range("H27").Formula = "=(H24-J24)/H23"
$ws1.Range("J27").Formula = "=SUM(B27:H27)/7"

# ---------------------------------------------------------------------
# 6) Duplicate the "r = (w - o)/ref w" note (row 11) to row 31.
# ---------------------------------------------------------------------
$ws1.Range("A31").Value = "r = (w - o)/ref w"

# ---------------------------------------------------------------------
# 7) Duplicate the second "Ref weights / Offsets" table (rows 14-15) to
#    rows 34-35.
# ---------------------------------------------------------------------
$ws1.Range("A34").Value = "Ref weights:"
$ws1.Range("B34").Value = 25
$ws1.Range("C34").Value = 100
$ws1.Range("D34").Value = 175
$ws1.Range("E34").Value = 250
$ws1.Range("J34").Value = "Offset (Sensor):"

$ws1.Range("A35").Value = "Offsets:"
$ws1.Range("B35").Value = 8671186.0600000005
$ws1.Range("C35").Value = 8794613.3066700008
$ws1.Range("D35").Value = 8918932.9733300004
$ws1.Range("E35").Value = 9043501.5133299995
$ws1.Range("J35").Value = 8628906.7799999993
$ws1.Range("B35:E35").NumberFormat = "0.000"
$ws1.Range("J35").NumberFormat = "0.000"

# ---------------------------------------------------------------------
# 8) Duplicate the "Average ratios / Ratios" table (rows 17-18) to rows
#    37-38.
# ---------------------------------------------------------------------
$ws1.Range("J37").Value = "Average ratios:"

$ws1.Range("A38").Value = "Ratios:"
$ws1.Range("B38").Formula = "=(B35-J35)/B34"
$ws1.Range("C38").Formula = "=(C35-J35)/C34"
$ws1.Range("D38").Formula = "=(D35-J35)/D34"
$ws1.Range("E38").Formula = "=(E35-J35)/E34"
$ws1.Range("J38").Formula = "=SUM(B38:E38)/4"

# ---------------------------------------------------------------------
# 9) View state: "1 Sensor" becomes the active tab, scrolled so row 4 is
#    the top-visible row, with P26 selected. "2 Sensors" loses the
#    active-tab flag it used to carry (but keeps its own B33 selection).
# ---------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("P26").Select()
